$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 373, shifting existing rows 373:416 down to 374:417
$ws.Rows.Item(373).Insert()

# Populate the newly inserted row 373 with the new data record
$ws.Cells.Item(373, 1).Value = 6
$ws.Cells.Item(373, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(373, 3).Value = "Metropolitana"
$ws.Cells.Item(373, 4).Value = 45127
$ws.Cells.Item(373, 5).Value = 13
$ws.Cells.Item(373, 6).Value = 100112026
$ws.Cells.Item(373, 7).Value = "Haba"
$ws.Cells.Item(373, 8).Value = "Sin especificar"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 650
$ws.Cells.Item(373, 11).Value = 15000
$ws.Cells.Item(373, 12).Value = 16000
$ws.Cells.Item(373, 13).Value = 15354
$ws.Cells.Item(373, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(373, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(373, 16).Value = 614
$ws.Cells.Item(373, 17).Value = 25
$ws.Cells.Item(373, 18).Value = "Hortaliza"
